$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'28.905.86"
$ws.Range("E2").Value = "  -2.90%  "

# Row 3
$ws.Range("D3").Value = "'1.881.90"
$ws.Range("E3").Value = "  -3.41%  "

# Row 4
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").Value = "'329.45"
$ws.Range("E5").Value = "  -3.72%  "

# Row 6
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.03%  "

# Row 7
$ws.Range("D7").Value = "'0.4583"
$ws.Range("E7").Value = "  -4.17%  "

# Row 8
$ws.Range("D8").Value = "'0.4105"
$ws.Range("E8").Value = "  -0.69%  "

# Row 9
$ws.Range("D9").Value = "'47.76"
$ws.Range("E9").Value = "  -2.29%  "

# Row 10
$ws.Range("D10").Value = "'0.07949"
$ws.Range("E10").Value = "  -3.81%  "

# Row 11
$ws.Range("D11").Value = "'0.9924"
$ws.Range("E11").Value = "  -5.18%  "

# Row 12
$ws.Range("D12").Value = "'21.60"
$ws.Range("E12").Value = "  -4.93%  "

# Row 13
$ws.Range("D13").Value = "'1.890.32"
$ws.Range("E13").Value = "  -2.47%  "

# Row 14
$ws.Range("D14").Value = "'5.902"
$ws.Range("E14").Value = "  -4.34%  "

# Row 15
$ws.Range("D15").Value = "'7.053"
$ws.Range("E15").Value = "  -5.23%  "

# Row 16
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  +0.00%  "

# Row 17
$ws.Range("D17").Value = "'88.37"
$ws.Range("E17").Value = "  -4.75%  "

# Row 18
$ws.Range("D18").Value = "'0.06574"
$ws.Range("E18").Value = "  -2.04%  "

# Row 19
$ws.Range("D19").Value = "'0.00001025"
$ws.Range("E19").Value = "  -3.79%  "

# Row 20
$ws.Range("D20").Value = "'17.37"
$ws.Range("E20").Value = "  -3.79%  "

# Row 21
$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "  +0.17%  "

# Row 22
$ws.Range("D22").Value = "'28.903.59"
$ws.Range("E22").Value = "  -2.77%  "

# Row 23
$ws.Range("D23").Value = "'5.406"
$ws.Range("E23").Value = "  -3.86%  "

# Row 24
$ws.Range("D24").Value = "'11.42"
$ws.Range("E24").Value = "  +1.28%  "

# Row 25
$ws.Range("E25").Value = "  -3.57%  "

# Row 26
$ws.Range("D26").Value = "'2.134.79"
$ws.Range("E26").Value = "  -1.58%  "

# Row 27
$ws.Range("D27").Value = "'156.04"
$ws.Range("E27").Value = "  -3.50%  "

# Row 28
$ws.Range("D28").Value = "'19.54"
$ws.Range("E28").Value = "  -3.24%  "

# Row 29
$ws.Range("D29").Value = "'2.079"
$ws.Range("E29").Value = "  -5.38%  "

# Row 30
$ws.Range("D30").Value = "'5.465"
$ws.Range("E30").Value = "  -3.23%  "

# Row 31
$ws.Range("D31").Value = "'117.41"
$ws.Range("E31").Value = "  -4.22%  "

# Row 32
$ws.Range("D32").Value = "'1.029"
$ws.Range("E32").Value = "  +0.26%  "

# Row 33
$ws.Range("D33").Value = "'0.09315"
$ws.Range("E33").Value = "  -3.47%  "

# Row 34
$ws.Range("D34").Value = "'1.398"
$ws.Range("E34").Value = "  -5.11%  "

# Row 35
$ws.Range("D35").Value = "'3.520"
$ws.Range("E35").Value = "  -4.40%  "

# Row 36
$ws.Range("D36").Value = "'5.281"
$ws.Range("E36").Value = "  -3.85%  "

# Row 37
$ws.Range("D37").Value = "'0.06055"
$ws.Range("E37").Value = "  -3.07%  "

# Row 38
$ws.Range("D38").Value = "'0.02224"
$ws.Range("E38").Value = "  -4.10%  "

# Row 39
$ws.Range("D39").Value = "'8.322"
$ws.Range("E39").Value = "  -4.51%  "

# Row 40
$ws.Range("D40").Value = "'1.171"
$ws.Range("E40").Value = "  -2.24%  "

# Row 41
$ws.Range("D41").Value = "'1.001"
$ws.Range("E41").Value = "  +0.03%  "

# Row 42
$ws.Range("D42").Value = "'0.5778"
$ws.Range("E42").Value = "  -5.41%  "

# Row 43
$ws.Range("D43").Value = "'0.1820"
$ws.Range("E43").Value = "  -4.60%  "

# Row 44
$ws.Range("D44").Value = "'10.05"
$ws.Range("E44").Value = "  -6.51%  "

# Row 45
$ws.Range("E45").Value = "  -2.63%  "

# Row 46
$ws.Range("D46").Value = "'0.07510"
$ws.Range("E46").Value = "  +0.80%  "

# Row 47
$ws.Range("D47").Value = "'2.261"
$ws.Range("E47").Value = "  -2.63%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'12.03"
$ws.Range("E48").Value = "  -4.10%  "

# Row 49
$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").Value = "'0.5443"
$ws.Range("E49").Value = "  -4.81%  "

# Row 50
$ws.Range("D50").Value = "'1.896"
$ws.Range("E50").Value = "  -5.18%  "

# Row 51
$ws.Range("D51").Value = "'111.18"
$ws.Range("E51").Value = "  -2.51%  "
